# Data Driving InvalidLogin Test case
# Adds a new "InvalidLogin" worksheet (after the existing "ValidLogin" sheet)
# populated with UserName/Password header row and a sample invalid
# credential pair, and refreshes the selections/active-sheet state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new sheet as the last tab in the workbook -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "InvalidLogin"

# --- Populate the InvalidLogin sheet with its header + sample data -------
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "ABC"
$ws2.Range("B2").Value = "XYZ"

# --- Refresh the selection on the original ValidLogin sheet --------------
$ws1.Activate()
$ws1.Range("A1:B2").Select()

# --- Make InvalidLogin the active tab with its own selection -------------
$ws2.Activate()
$ws2.Range("O12").Select()
